$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

# Row 13: F13 effort hrs 0.41 -> 1.41
$ws.Cells.Item(13, 6).Value = 1.41

# Row 21: F21 effort hrs 0 -> 1
$ws.Cells.Item(21, 6).Value = 1

# Row 22: F22 effort hrs 0 -> 1
$ws.Cells.Item(22, 6).Value = 1

# Row 23 (CU-19 "CONSULTAR EGRESOS"): mark as discarded and explain why
$ws.Cells.Item(23, 5).Value = "DESCARTADO"
$ws.Cells.Item(23, 9).Value = "El CU fue descartado ya que se repite en el CU 20 Generar reportes mensuales"

# Update the visible selection to match the saved view (F21)
$ws.Range("F21").Select()
